$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before column A (Excel's "Insert Sheet Columns"),
# shifting all existing data one column to the right.
$ws.Columns("A").Insert()

# Populate the new column A: a header label plus a per-row sequence number
# for each of the 8 data rows.
$ws.Range("A1").Value = "population_density"
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 5
$ws.Range("A7").Value = 6
$ws.Range("A8").Value = 7
$ws.Range("A9").Value = 8

# Restore the selection to where the author left it.
[void]$ws.Range("A10").Select()
